$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Force the Price/Volume columns to Text format before writing so Excel
# does not reinterpret numeric-looking strings (e.g. "1.00", "65.544.67")
# as numbers, matching the source inlineStr cells exactly.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '65.544.67'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '2.659.32'
$ws.Range('E3').Value = '  -0.92%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '596.48'
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('D6').Value = '155.75'
$ws.Range('E6').Value = '  -1.10%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.619'
$ws.Range('E8').Value = '  +5.35%  '
$ws.Range('D9').Value = '0.127'
$ws.Range('E9').Value = '  +2.80%  '
$ws.Range('D10').Value = '0.396'
$ws.Range('E10').Value = '  -0.87%  '
$ws.Range('D11').Value = '5.79'
$ws.Range('E11').Value = '  -2.59%  '
$ws.Range('E12').Value = '  +0.16%  '
$ws.Range('D13').Value = '28.81'
$ws.Range('E13').Value = '  -3.33%  '
$ws.Range('D14').Value = '0.0000196'
$ws.Range('E14').Value = '  -4.12%  '
$ws.Range('D15').Value = '3.138.81'
$ws.Range('E15').Value = '  -0.82%  '
$ws.Range('D16').Value = '65.435.47'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').Value = '2.673.34'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').Value = '12.61'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('D19').Value = '4.76'
$ws.Range('E19').Value = '  -1.76%  '
$ws.Range('D20').Value = '7.46'
$ws.Range('E20').Value = '  -2.63%  '
$ws.Range('D21').Value = '348.73'
$ws.Range('E21').Value = '  -1.13%  '
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').Value = '68.95'
$ws.Range('E23').Value = '  -2.18%  '
$ws.Range('D24').Value = '0.0000112'
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('D25').Value = '9.65'
$ws.Range('E25').Value = '  -1.95%  '
$ws.Range('E26').Value = '  +2.61%  '
$ws.Range('D27').Value = '1.58'
$ws.Range('E27').Value = '  -3.33%  '
$ws.Range('E28').Value = '  -3.74%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').Value = '7.95'
$ws.Range('E30').Value = '  -3.45%  '
$ws.Range('D31').Value = '2.12'
$ws.Range('E31').Value = '  -2.83%  '
$ws.Range('D32').Value = '530.30'
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('D33').Value = '1.76'
$ws.Range('E33').Value = '  -0.79%  '
$ws.Range('D34').Value = '6.40'
$ws.Range('E34').Value = '  -3.25%  '
$ws.Range('D35').Value = '5.41'
$ws.Range('E35').Value = '  -0.51%  '
$ws.Range('D36').Value = '0.420'
$ws.Range('E36').Value = '  -2.25%  '
$ws.Range('D37').Value = '20.42'
$ws.Range('E37').Value = '  -0.89%  '
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('D39').Value = '155.38'
$ws.Range('E39').Value = '  -3.67%  '
$ws.Range('D40').Value = '1.92'
$ws.Range('E40').Value = '  -3.18%  '
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').Value = '161.13'
$ws.Range('E42').Value = '  -3.37%  '
$ws.Range('D43').Value = '4.07'
$ws.Range('E43').Value = '  -1.25%  '
$ws.Range('D44').Value = '2.31'
$ws.Range('E44').Value = '  +2.50%  '
$ws.Range('D45').Value = '0.0606'
$ws.Range('E45').Value = '  -3.48%  '
$ws.Range('D46').Value = '22.50'
$ws.Range('E46').Value = '  -3.91%  '
$ws.Range('D47').Value = '0.636'
$ws.Range('E47').Value = '  -2.47%  '
$ws.Range('D48').Value = '0.0256'
$ws.Range('E48').Value = '  -3.22%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '0.0987'
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0250'
$ws.Range('E50').Value = '  +5.50%  '
$ws.Range('D51').Value = '19.76'
$ws.Range('E51').Value = '  -3.81%  '

# Restore the default (unstyled) cell style now that the text values are set,
# so no stray formatting is left behind on these data cells.
$ws.Range("D2:E51").Style = "Normal"
